# Applies the coinranking.com crypto price/volume refresh described in the
# commit "Updated cryptos list ... with GitHub Actions".
# Row 30/31 additionally swap places (PancakeSwap <-> NEARProtocol).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "67.726.32"
$ws.Range("E2").Value = "  +4.49%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "3.267.47"
$ws.Range("E3").Value = "  +4.36%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  +0.04%  "

# Row 5: BNB
$ws.Range("D5").Value = "'580.19"
$ws.Range("E5").Value = "  +2.25%  "

# Row 6: Solana
$ws.Range("E6").Value = "  +8.69%  "

# Row 7: USDC
$ws.Range("E7").Value = "  +0.20%  "

# Row 8: XRP
$ws.Range("E8").Value = "  -0.26%  "

# Row 9: LidoStakedEther
$ws.Range("D9").Value = "3.266.77"
$ws.Range("E9").Value = "  +4.53%  "

# Row 10: Dogecoin
$ws.Range("E10").Value = "  +8.27%  "

# Row 11: Toncoin
$ws.Range("D11").Value = "'6.75"
$ws.Range("E11").Value = "  +3.67%  "

# Row 12: Cardano
$ws.Range("E12").Value = "  +7.56%  "

# Row 13: WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "3.834.25"
$ws.Range("E13").Value = "  +4.86%  "

# Row 14: TRON
$ws.Range("E14").Value = "  +1.38%  "

# Row 15: Avalanche
$ws.Range("D15").Value = "'28.63"
$ws.Range("E15").Value = "  +7.73%  "

# Row 16: WrappedBTC
$ws.Range("D16").Value = "67.705.12"
$ws.Range("E16").Value = "  +4.63%  "

# Row 17: ShibaInu
$ws.Range("D17").Value = "'0.0000168"
$ws.Range("E17").Value = "  +4.90%  "

# Row 18: WrappedEther
$ws.Range("D18").Value = "3.270.51"
$ws.Range("E18").Value = "  +4.85%  "

# Row 19: Polkadot
$ws.Range("D19").Value = "'5.85"
$ws.Range("E19").Value = "  +3.31%  "

# Row 20: Chainlink
$ws.Range("E20").Value = "  +7.71%  "

# Row 21: BitcoinCash
$ws.Range("D21").Value = "'374.93"
$ws.Range("E21").Value = "  +5.92%  "

# Row 22: Uniswap
$ws.Range("D22").Value = "'7.66"
$ws.Range("E22").Value = "  +6.50%  "

# Row 23: Dai
$ws.Range("E23").Value = "  +0.07%  "

# Row 24: Litecoin
$ws.Range("D24").Value = "'71.28"
$ws.Range("E24").Value = "  +3.90%  "

# Row 25: Polygon
$ws.Range("D25").Value = "'0.513"
$ws.Range("E25").Value = "  +4.51%  "

# Row 26: PEPE
$ws.Range("E26").Value = "  +5.84%  "

# Row 27: InternetComputer(DFINITY)
$ws.Range("D27").Value = "'9.66"
$ws.Range("E27").Value = "  +0.35%  "

# Row 28: Kaspa
$ws.Range("E28").Value = "  +3.67%  "

# Row 29: Binance-PegBSC-USD
$ws.Range("E29").Value = "  +0.08%  "

# Row 30: NEARProtocol
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'1.98"
$ws.Range("E30").Value = "  +4.38%  "

# Row 31: PancakeSwap
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "'5.72"
$ws.Range("E31").Value = "  +9.48%  "

# Row 32: EthereumClassic
$ws.Range("D32").Value = "'22.78"
$ws.Range("E32").Value = "  +5.02%  "

# Row 33: USDe
$ws.Range("E33").Value = "  +0.00%  "

# Row 34: Fetch.AI
$ws.Range("E34").Value = "  +8.29%  "

# Row 35: Aptos
$ws.Range("D35").Value = "'6.94"
$ws.Range("E35").Value = "  +6.18%  "

# Row 36: ImmutableX
$ws.Range("E36").Value = "  +6.29%  "

# Row 37: Monero
$ws.Range("E37").Value = "  +3.21%  "

# Row 38: Mantle
$ws.Range("D38").Value = "'0.853"
$ws.Range("E38").Value = "  +3.62%  "

# Row 39: Stacks
$ws.Range("E39").Value = "  +5.92%  "

# Row 40: RenderToken
$ws.Range("D40").Value = "'6.85"
$ws.Range("E40").Value = "  +12.79%  "

# Row 41: Filecoin
$ws.Range("E41").Value = "  +13.53%  "

# Row 42: EnergySwap
$ws.Range("D42").Value = "'26.86"
$ws.Range("E42").Value = "  +4.70%  "

# Row 43: dogwifhat
$ws.Range("E43").Value = "  +10.14%  "

# Row 44: Bittensor
$ws.Range("D44").Value = "'358.62"
$ws.Range("E44").Value = "  +13.12%  "

# Row 45: Maker
$ws.Range("D45").Value = "2.711.08"
$ws.Range("E45").Value = "  +2.90%  "

# Row 46: InjectiveProtocol
$ws.Range("D46").Value = "'25.51"
$ws.Range("E46").Value = "  +7.33%  "

# Row 47: OKB
$ws.Range("E47").Value = "  +4.24%  "

# Row 48: Hedera
$ws.Range("E48").Value = "  +5.68%  "

# Row 49: VeChain
$ws.Range("D49").Value = "'0.0282"
$ws.Range("E49").Value = "  +4.54%  "

# Row 50: ONDO
$ws.Range("D50").Value = "'1.01"
$ws.Range("E50").Value = "  +8.31%  "

# Row 51: Stellar
$ws.Range("E51").Value = "  +0.91%  "

Write-Host "cryptos list updated"
